$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (column F) counts
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 5560
$wsExhibit.Range("F10").Value = 2490
$wsExhibit.Range("F12").Value = 110
$wsExhibit.Range("F15").Value = 8
$wsExhibit.Range("F17").Value = 280

# Sheet "全部类型" (sheet4): same events mirrored, update matching counts
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 5560
$wsAll.Range("F12").Value = 2490
$wsAll.Range("F14").Value = 110
$wsAll.Range("F18").Value = 8
$wsAll.Range("F20").Value = 280
